$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 11
$srcRow = 10

# Copy the formatting of the previous year's row down onto the new row
# first (so the new "year" label cell picks up the same centered/bold/
# bordered style as A2:A10), then fill in the actual values.
$ws.Range($ws.Cells.Item($srcRow, 1), $ws.Cells.Item($srcRow, 43)).Copy() | Out-Null
$ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 43)).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = "2021年"

$values = @{
    2  = 1201.38
    3  = 286.55
    4  = 53.86
    6  = 1363.35
    7  = 1381.61
    8  = 284.85
    9  = 713.95
    10 = 162.15
    11 = 250.16
    12 = 230.34
    13 = 1.61
    14 = 511.67
    15 = 835.67
    16 = 53.15
    17 = 348.12
    18 = 856.48
    19 = 5.64
    20 = 997.63
    21 = 4.79
    22 = 265.82
    23 = 41.23
    24 = 21.18
    25 = 2302.71
    26 = 257.18
    27 = 670.9299999999999
    28 = 0.47
    29 = 23935.16
    30 = 1128.3
    31 = 515.77
    32 = 1684.14
    33 = 1406.68
    34 = 294.35
    35 = 416.44
    36 = 6.47
    37 = 1355.83
    38 = 252.99
    39 = 1761.21
    40 = 97.09999999999999
    41 = 415.27
    42 = 1372.19
    43 = 125.96
}

foreach ($col in $values.Keys) {
    $ws.Cells.Item($row, $col).Value = $values[$col]
}

# Column E ("其他采矿业私营工业企业产成品") has no reported figure for 2021,
# same as the prior year (E10) — leave it blank.
